$d = $word.ActiveDocument

$replacements = @(
    @("2025-08-15 Friday", "2025-08-16 Saturday"),
    @("31×19=", "27×65="),
    @("24×19=", "72×39="),
    @("31×17=", "14×11="),
    @("71×99=", "54×90="),
    @("88×24=", "60×55="),
    @("20×62=", "46×31="),
    @("22×45=", "63×41="),
    @("72×19=", "78×74="),
    @("70×58=", "64×46="),
    @("60×12=", "73×97="),
    @("88×99=", "51×44="),
    @("25×69=", "24×37="),
    @("37×37=", "46×17="),
    @("29×66=", "57×83="),
    @("77×62=", "12×69="),
    @("21×67=", "22×70="),
    @("33×72=", "15×87="),
    @("82×30=", "23×70="),
    @("59×34=", "73×21="),
    @("95×83=", "33×11="),
    @("38×81=", "88×56="),
    @("18×11=", "15×61="),
    @("59×79=", "54×64="),
    @("92×12=", "94×54="),
    @("30×55=", "70×83="),
)

foreach ($pair in $replacements) {
    $oldText = $pair[0]
    $newText = $pair[1]
    $ok = $d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
    if (-not $ok) {
        throw "Failed to replace: $oldText -> $newText"
    }
    Write-Output "Replaced $oldText -> $newText"
}
